$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained a new "2020" column (Q), to the right of the existing
# "2019" column (P). Copy column P's cell formatting into column Q first
# (so Q4/Q5/Q6/Q7/Q8 pick up the same styles as P4/P5/P6/P7/P8), then set
# the new column's values.
$ws.Range("P4:P8").Copy($ws.Range("Q4:Q8"))

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 2
$ws.Range("Q6").Value = 0.3
$ws.Range("Q7").Value = 0.1
$ws.Range("Q8").Value = 4.3

# Move the active selection, matching the saved view state.
$ws.Range("O12").Select()
